# Weekly data refresh for the Repollo price table.
# Existing records in rows 88-214 shift down by one row (each row now
# holds the record that used to be one row below it); a brand-new record
# is written into row 88, and the record that used to be last (old row
# 214) becomes the new last row, 215.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A88:A215").Value = 5
$ws.Range("B88:B215").Value = 'Macroferia Regional de Talca'
$ws.Range("C88:C215").Value = 'Maule'
$ws.Range("E88:E215").Value = 7
$ws.Range("F88:F215").Value = 100112006
$ws.Range("G88:G215").Value = 'Repollo'
$ws.Range("N88:N215").Value = '$/unidad'
$ws.Range("Q88:Q215").Value = 1
$ws.Range("R88:R215").Value = 'Hortaliza'

$D = New-Object 'object[,]' 128,1
$D[0,0] = 44557
$D[1,0] = 44384
$D[2,0] = 44320
$D[3,0] = 44294
$D[4,0] = 44326
$D[5,0] = 44354
$D[6,0] = 44315
$D[7,0] = 44505
$D[8,0] = 44274
$D[9,0] = 44348
$D[10,0] = 44309
$D[11,0] = 44214
$D[12,0] = 44172
$D[13,0] = 44546
$D[14,0] = 44312
$D[15,0] = 44349
$D[16,0] = 44376
$D[17,0] = 44166
$D[18,0] = 44327
$D[19,0] = 44316
$D[20,0] = 44277
$D[21,0] = 44390
$D[22,0] = 44482
$D[23,0] = 44273
$D[24,0] = 44364
$D[25,0] = 44469
$D[26,0] = 44265
$D[27,0] = 44330
$D[28,0] = 44280
$D[29,0] = 44280
$D[30,0] = 44323
$D[31,0] = 44323
$D[32,0] = 44525
$D[33,0] = 44266
$D[34,0] = 44403
$D[35,0] = 44403
$D[36,0] = 44165
$D[37,0] = 44267
$D[38,0] = 44306
$D[39,0] = 44516
$D[40,0] = 44427
$D[41,0] = 44495
$D[42,0] = 44176
$D[43,0] = 44370
$D[44,0] = 44475
$D[45,0] = 44508
$D[46,0] = 44291
$D[47,0] = 44468
$D[48,0] = 44379
$D[49,0] = 44498
$D[50,0] = 44449
$D[51,0] = 44449
$D[52,0] = 44526
$D[53,0] = 44271
$D[54,0] = 44425
$D[55,0] = 44425
$D[56,0] = 44168
$D[57,0] = 44447
$D[58,0] = 44553
$D[59,0] = 44167
$D[60,0] = 44161
$D[61,0] = 44328
$D[62,0] = 44328
$D[63,0] = 44285
$D[64,0] = 44418
$D[65,0] = 44418
$D[66,0] = 44160
$D[67,0] = 44434
$D[68,0] = 44434
$D[69,0] = 44467
$D[70,0] = 44231
$D[71,0] = 44490
$D[72,0] = 44341
$D[73,0] = 44279
$D[74,0] = 44264
$D[75,0] = 44322
$D[76,0] = 44322
$D[77,0] = 44491
$D[78,0] = 44389
$D[79,0] = 44391
$D[80,0] = 44396
$D[81,0] = 44510
$D[82,0] = 44551
$D[83,0] = 44386
$D[84,0] = 44519
$D[85,0] = 44420
$D[86,0] = 44420
$D[87,0] = 44414
$D[88,0] = 44414
$D[89,0] = 44543
$D[90,0] = 44321
$D[91,0] = 44278
$D[92,0] = 44308
$D[93,0] = 44281
$D[94,0] = 44281
$D[95,0] = 44474
$D[96,0] = 44446
$D[97,0] = 44446
$D[98,0] = 44350
$D[99,0] = 44529
$D[100,0] = 44405
$D[101,0] = 44405
$D[102,0] = 44413
$D[103,0] = 44413
$D[104,0] = 44411
$D[105,0] = 44411
$D[106,0] = 44258
$D[107,0] = 44298
$D[108,0] = 44298
$D[109,0] = 44432
$D[110,0] = 44428
$D[111,0] = 44340
$D[112,0] = 44224
$D[113,0] = 44329
$D[114,0] = 44452
$D[115,0] = 44452
$D[116,0] = 44270
$D[117,0] = 44398
$D[118,0] = 44536
$D[119,0] = 44302
$D[120,0] = 44344
$D[121,0] = 44344
$D[122,0] = 44461
$D[123,0] = 44357
$D[124,0] = 44371
$D[125,0] = 44365
$D[126,0] = 44313
$D[127,0] = 44272
$ws.Range("D88:D215").Value = $D

$H = New-Object 'object[,]' 128,1
$H[0,0] = 'Crespo record'
$H[1,0] = 'Crespo record'
$H[2,0] = 'Crespo record'
$H[3,0] = 'Crespo record'
$H[4,0] = 'Crespo record'
$H[5,0] = 'Crespo record'
$H[6,0] = 'Crespo record'
$H[7,0] = 'Crespo record'
$H[8,0] = 'Crespo record'
$H[9,0] = 'Crespo record'
$H[10,0] = 'Crespo record'
$H[11,0] = 'Crespo record'
$H[12,0] = 'Crespo record'
$H[13,0] = 'Crespo record'
$H[14,0] = 'Crespo record'
$H[15,0] = 'Crespo record'
$H[16,0] = 'Crespo record'
$H[17,0] = 'Crespo record'
$H[18,0] = 'Crespo record'
$H[19,0] = 'Crespo record'
$H[20,0] = 'Crespo record'
$H[21,0] = 'Crespo record'
$H[22,0] = 'Crespo record'
$H[23,0] = 'Crespo record'
$H[24,0] = 'Crespo record'
$H[25,0] = 'Crespo record'
$H[26,0] = 'Crespo record'
$H[27,0] = 'Crespo record'
$H[28,0] = 'Crespo record'
$H[29,0] = 'Crespo record'
$H[30,0] = 'Crespo record'
$H[31,0] = 'Crespo record'
$H[32,0] = 'Crespo record'
$H[33,0] = 'Crespo record'
$H[34,0] = 'Crespo record'
$H[35,0] = 'Crespo record'
$H[36,0] = 'Crespo record'
$H[37,0] = 'Crespo record'
$H[38,0] = 'Crespo record'
$H[39,0] = 'Crespo record'
$H[40,0] = 'Crespo record'
$H[41,0] = 'Crespo record'
$H[42,0] = 'Crespo record'
$H[43,0] = 'Crespo record'
$H[44,0] = 'Crespo record'
$H[45,0] = 'Crespo record'
$H[46,0] = 'Crespo record'
$H[47,0] = 'Crespo record'
$H[48,0] = 'Crespo record'
$H[49,0] = 'Crespo record'
$H[50,0] = 'Crespo record'
$H[51,0] = 'Crespo record'
$H[52,0] = 'Crespo record'
$H[53,0] = 'Crespo record'
$H[54,0] = 'Crespo record'
$H[55,0] = 'Crespo record'
$H[56,0] = 'Crespo record'
$H[57,0] = 'Crespo record'
$H[58,0] = 'Crespo record'
$H[59,0] = 'Crespo record'
$H[60,0] = 'Crespo record'
$H[61,0] = 'Crespo record'
$H[62,0] = 'Crespo record'
$H[63,0] = 'Crespo record'
$H[64,0] = 'Crespo record'
$H[65,0] = 'Crespo record'
$H[66,0] = 'Crespo record'
$H[67,0] = 'Crespo record'
$H[68,0] = 'Crespo record'
$H[69,0] = 'Crespo record'
$H[70,0] = 'Crespo record'
$H[71,0] = 'Crespo record'
$H[72,0] = 'Crespo record'
$H[73,0] = 'Crespo record'
$H[74,0] = 'Crespo record'
$H[75,0] = 'Crespo record'
$H[76,0] = 'Crespo record'
$H[77,0] = 'Crespo record'
$H[78,0] = 'Crespo record'
$H[79,0] = 'Crespo record'
$H[80,0] = 'Crespo record'
$H[81,0] = 'Crespo record'
$H[82,0] = 'Crespo record'
$H[83,0] = 'Crespo record'
$H[84,0] = 'Crespo record'
$H[85,0] = 'Crespo record'
$H[86,0] = 'Crespo record'
$H[87,0] = 'Crespo record'
$H[88,0] = 'Crespo record'
$H[89,0] = 'Crespo record'
$H[90,0] = 'Crespo record'
$H[91,0] = 'Crespo record'
$H[92,0] = 'Crespo record'
$H[93,0] = 'Crespo record'
$H[94,0] = 'Crespo record'
$H[95,0] = 'Crespo record'
$H[96,0] = 'Crespo record'
$H[97,0] = 'Crespo record'
$H[98,0] = 'Crespo record'
$H[99,0] = 'Crespo record'
$H[100,0] = 'Crespo record'
$H[101,0] = 'Crespo record'
$H[102,0] = 'Crespo record'
$H[103,0] = 'Crespo record'
$H[104,0] = 'Crespo record'
$H[105,0] = 'Crespo record'
$H[106,0] = 'Copenhague'
$H[107,0] = 'Crespo record'
$H[108,0] = 'Crespo record'
$H[109,0] = 'Crespo record'
$H[110,0] = 'Crespo record'
$H[111,0] = 'Crespo record'
$H[112,0] = 'Crespo record'
$H[113,0] = 'Crespo record'
$H[114,0] = 'Crespo record'
$H[115,0] = 'Crespo record'
$H[116,0] = 'Crespo record'
$H[117,0] = 'Crespo record'
$H[118,0] = 'Crespo record'
$H[119,0] = 'Crespo record'
$H[120,0] = 'Crespo record'
$H[121,0] = 'Crespo record'
$H[122,0] = 'Crespo record'
$H[123,0] = 'Crespo record'
$H[124,0] = 'Crespo record'
$H[125,0] = 'Crespo record'
$H[126,0] = 'Crespo record'
$H[127,0] = 'Crespo record'
$ws.Range("H88:H215").Value = $H

$I = New-Object 'object[,]' 128,1
$I[0,0] = 'Primera'
$I[1,0] = 'Primera'
$I[2,0] = 'Primera'
$I[3,0] = 'Primera'
$I[4,0] = 'Primera'
$I[5,0] = 'Primera'
$I[6,0] = 'Primera'
$I[7,0] = 'Primera'
$I[8,0] = 'Primera'
$I[9,0] = 'Primera'
$I[10,0] = 'Primera'
$I[11,0] = 'Primera'
$I[12,0] = 'Primera'
$I[13,0] = 'Primera'
$I[14,0] = 'Primera'
$I[15,0] = 'Primera'
$I[16,0] = 'Primera'
$I[17,0] = 'Primera'
$I[18,0] = 'Primera'
$I[19,0] = 'Primera'
$I[20,0] = 'Primera'
$I[21,0] = 'Primera'
$I[22,0] = 'Primera'
$I[23,0] = 'Primera'
$I[24,0] = 'Primera'
$I[25,0] = 'Primera'
$I[26,0] = 'Primera'
$I[27,0] = 'Primera'
$I[28,0] = 'Primera'
$I[29,0] = 'Segunda'
$I[30,0] = 'Primera'
$I[31,0] = 'Segunda'
$I[32,0] = 'Primera'
$I[33,0] = 'Primera'
$I[34,0] = 'Primera'
$I[35,0] = 'Segunda'
$I[36,0] = 'Primera'
$I[37,0] = 'Primera'
$I[38,0] = 'Primera'
$I[39,0] = 'Primera'
$I[40,0] = 'Segunda'
$I[41,0] = 'Primera'
$I[42,0] = 'Primera'
$I[43,0] = 'Primera'
$I[44,0] = 'Primera'
$I[45,0] = 'Primera'
$I[46,0] = 'Primera'
$I[47,0] = 'Primera'
$I[48,0] = 'Primera'
$I[49,0] = 'Primera'
$I[50,0] = 'Primera'
$I[51,0] = 'Segunda'
$I[52,0] = 'Primera'
$I[53,0] = 'Primera'
$I[54,0] = 'Primera'
$I[55,0] = 'Segunda'
$I[56,0] = 'Primera'
$I[57,0] = 'Primera'
$I[58,0] = 'Primera'
$I[59,0] = 'Primera'
$I[60,0] = 'Primera'
$I[61,0] = 'Primera'
$I[62,0] = 'Segunda'
$I[63,0] = 'Primera'
$I[64,0] = 'Primera'
$I[65,0] = 'Segunda'
$I[66,0] = 'Primera'
$I[67,0] = 'Primera'
$I[68,0] = 'Segunda'
$I[69,0] = 'Primera'
$I[70,0] = 'Primera'
$I[71,0] = 'Primera'
$I[72,0] = 'Primera'
$I[73,0] = 'Primera'
$I[74,0] = 'Primera'
$I[75,0] = 'Primera'
$I[76,0] = 'Segunda'
$I[77,0] = 'Primera'
$I[78,0] = 'Primera'
$I[79,0] = 'Primera'
$I[80,0] = 'Primera'
$I[81,0] = 'Primera'
$I[82,0] = 'Primera'
$I[83,0] = 'Primera'
$I[84,0] = 'Primera'
$I[85,0] = 'Primera'
$I[86,0] = 'Segunda'
$I[87,0] = 'Primera'
$I[88,0] = 'Segunda'
$I[89,0] = 'Primera'
$I[90,0] = 'Primera'
$I[91,0] = 'Primera'
$I[92,0] = 'Primera'
$I[93,0] = 'Primera'
$I[94,0] = 'Segunda'
$I[95,0] = 'Primera'
$I[96,0] = 'Primera'
$I[97,0] = 'Segunda'
$I[98,0] = 'Primera'
$I[99,0] = 'Primera'
$I[100,0] = 'Primera'
$I[101,0] = 'Segunda'
$I[102,0] = 'Primera'
$I[103,0] = 'Segunda'
$I[104,0] = 'Primera'
$I[105,0] = 'Segunda'
$I[106,0] = 'Primera'
$I[107,0] = 'Primera'
$I[108,0] = 'Segunda'
$I[109,0] = 'Segunda'
$I[110,0] = 'Segunda'
$I[111,0] = 'Primera'
$I[112,0] = 'Primera'
$I[113,0] = 'Primera'
$I[114,0] = 'Primera'
$I[115,0] = 'Segunda'
$I[116,0] = 'Primera'
$I[117,0] = 'Primera'
$I[118,0] = 'Primera'
$I[119,0] = 'Primera'
$I[120,0] = 'Primera'
$I[121,0] = 'Segunda'
$I[122,0] = 'Primera'
$I[123,0] = 'Primera'
$I[124,0] = 'Primera'
$I[125,0] = 'Primera'
$I[126,0] = 'Primera'
$I[127,0] = 'Primera'
$ws.Range("I88:I215").Value = $I

$J = New-Object 'object[,]' 128,1
$J[0,0] = 5000
$J[1,0] = 5000
$J[2,0] = 4000
$J[3,0] = 4000
$J[4,0] = 5000
$J[5,0] = 5000
$J[6,0] = 4000
$J[7,0] = 6000
$J[8,0] = 2000
$J[9,0] = 5000
$J[10,0] = 3000
$J[11,0] = 3000
$J[12,0] = 1500
$J[13,0] = 5000
$J[14,0] = 5000
$J[15,0] = 5000
$J[16,0] = 4000
$J[17,0] = 4000
$J[18,0] = 5000
$J[19,0] = 5000
$J[20,0] = 3000
$J[21,0] = 5000
$J[22,0] = 4000
$J[23,0] = 3000
$J[24,0] = 5000
$J[25,0] = 5000
$J[26,0] = 2000
$J[27,0] = 3000
$J[28,0] = 2000
$J[29,0] = 2000
$J[30,0] = 2000
$J[31,0] = 2000
$J[32,0] = 3000
$J[33,0] = 2000
$J[34,0] = 2000
$J[35,0] = 4000
$J[36,0] = 3000
$J[37,0] = 3000
$J[38,0] = 3000
$J[39,0] = 3000
$J[40,0] = 3000
$J[41,0] = 5000
$J[42,0] = 2000
$J[43,0] = 5000
$J[44,0] = 3000
$J[45,0] = 6000
$J[46,0] = 3000
$J[47,0] = 3000
$J[48,0] = 4000
$J[49,0] = 6000
$J[50,0] = 3000
$J[51,0] = 2000
$J[52,0] = 6000
$J[53,0] = 3000
$J[54,0] = 2000
$J[55,0] = 3000
$J[56,0] = 2000
$J[57,0] = 4000
$J[58,0] = 4000
$J[59,0] = 3000
$J[60,0] = 3000
$J[61,0] = 3000
$J[62,0] = 2000
$J[63,0] = 3000
$J[64,0] = 2000
$J[65,0] = 3000
$J[66,0] = 2000
$J[67,0] = 2000
$J[68,0] = 3000
$J[69,0] = 3000
$J[70,0] = 3000
$J[71,0] = 3000
$J[72,0] = 5000
$J[73,0] = 3000
$J[74,0] = 3000
$J[75,0] = 3000
$J[76,0] = 3000
$J[77,0] = 5000
$J[78,0] = 5000
$J[79,0] = 5000
$J[80,0] = 5000
$J[81,0] = 6000
$J[82,0] = 3000
$J[83,0] = 8000
$J[84,0] = 4000
$J[85,0] = 2000
$J[86,0] = 5000
$J[87,0] = 2000
$J[88,0] = 5000
$J[89,0] = 5000
$J[90,0] = 5000
$J[91,0] = 3000
$J[92,0] = 3000
$J[93,0] = 2000
$J[94,0] = 1000
$J[95,0] = 4000
$J[96,0] = 2000
$J[97,0] = 3000
$J[98,0] = 5000
$J[99,0] = 3000
$J[100,0] = 2000
$J[101,0] = 4000
$J[102,0] = 3000
$J[103,0] = 3000
$J[104,0] = 2000
$J[105,0] = 4000
$J[106,0] = 2000
$J[107,0] = 2000
$J[108,0] = 2000
$J[109,0] = 3000
$J[110,0] = 3000
$J[111,0] = 3000
$J[112,0] = 3000
$J[113,0] = 4000
$J[114,0] = 2000
$J[115,0] = 5000
$J[116,0] = 3000
$J[117,0] = 5000
$J[118,0] = 3000
$J[119,0] = 3000
$J[120,0] = 3000
$J[121,0] = 2000
$J[122,0] = 3000
$J[123,0] = 5000
$J[124,0] = 5000
$J[125,0] = 3000
$J[126,0] = 5000
$J[127,0] = 2000
$ws.Range("J88:J215").Value = $J

$K = New-Object 'object[,]' 128,1
$K[0,0] = 600
$K[1,0] = 500
$K[2,0] = 500
$K[3,0] = 800
$K[4,0] = 400
$K[5,0] = 500
$K[6,0] = 400
$K[7,0] = 600
$K[8,0] = 800
$K[9,0] = 500
$K[10,0] = 600
$K[11,0] = 900
$K[12,0] = 800
$K[13,0] = 500
$K[14,0] = 400
$K[15,0] = 500
$K[16,0] = 500
$K[17,0] = 700
$K[18,0] = 450
$K[19,0] = 400
$K[20,0] = 800
$K[21,0] = 400
$K[22,0] = 600
$K[23,0] = 800
$K[24,0] = 400
$K[25,0] = 400
$K[26,0] = 900
$K[27,0] = 500
$K[28,0] = 800
$K[29,0] = 600
$K[30,0] = 450
$K[31,0] = 300
$K[32,0] = 800
$K[33,0] = 800
$K[34,0] = 500
$K[35,0] = 350
$K[36,0] = 700
$K[37,0] = 800
$K[38,0] = 700
$K[39,0] = 900
$K[40,0] = 350
$K[41,0] = 600
$K[42,0] = 700
$K[43,0] = 400
$K[44,0] = 500
$K[45,0] = 700
$K[46,0] = 700
$K[47,0] = 450
$K[48,0] = 500
$K[49,0] = 600
$K[50,0] = 500
$K[51,0] = 300
$K[52,0] = 800
$K[53,0] = 800
$K[54,0] = 450
$K[55,0] = 300
$K[56,0] = 700
$K[57,0] = 500
$K[58,0] = 700
$K[59,0] = 700
$K[60,0] = 900
$K[61,0] = 500
$K[62,0] = 350
$K[63,0] = 800
$K[64,0] = 500
$K[65,0] = 350
$K[66,0] = 900
$K[67,0] = 500
$K[68,0] = 350
$K[69,0] = 500
$K[70,0] = 800
$K[71,0] = 700
$K[72,0] = 500
$K[73,0] = 800
$K[74,0] = 1000
$K[75,0] = 450
$K[76,0] = 350
$K[77,0] = 700
$K[78,0] = 400
$K[79,0] = 400
$K[80,0] = 350
$K[81,0] = 800
$K[82,0] = 900
$K[83,0] = 500
$K[84,0] = 900
$K[85,0] = 500
$K[86,0] = 350
$K[87,0] = 500
$K[88,0] = 350
$K[89,0] = 500
$K[90,0] = 400
$K[91,0] = 800
$K[92,0] = 600
$K[93,0] = 800
$K[94,0] = 600
$K[95,0] = 500
$K[96,0] = 500
$K[97,0] = 300
$K[98,0] = 500
$K[99,0] = 800
$K[100,0] = 500
$K[101,0] = 350
$K[102,0] = 500
$K[103,0] = 350
$K[104,0] = 500
$K[105,0] = 400
$K[106,0] = 800
$K[107,0] = 800
$K[108,0] = 600
$K[109,0] = 350
$K[110,0] = 300
$K[111,0] = 500
$K[112,0] = 800
$K[113,0] = 500
$K[114,0] = 500
$K[115,0] = 300
$K[116,0] = 800
$K[117,0] = 350
$K[118,0] = 800
$K[119,0] = 800
$K[120,0] = 600
$K[121,0] = 400
$K[122,0] = 500
$K[123,0] = 500
$K[124,0] = 450
$K[125,0] = 400
$K[126,0] = 400
$K[127,0] = 800
$ws.Range("K88:K215").Value = $K

$L = New-Object 'object[,]' 128,1
$L[0,0] = 600
$L[1,0] = 500
$L[2,0] = 500
$L[3,0] = 800
$L[4,0] = 400
$L[5,0] = 500
$L[6,0] = 400
$L[7,0] = 600
$L[8,0] = 800
$L[9,0] = 500
$L[10,0] = 600
$L[11,0] = 900
$L[12,0] = 800
$L[13,0] = 500
$L[14,0] = 400
$L[15,0] = 500
$L[16,0] = 500
$L[17,0] = 700
$L[18,0] = 450
$L[19,0] = 400
$L[20,0] = 800
$L[21,0] = 400
$L[22,0] = 600
$L[23,0] = 800
$L[24,0] = 400
$L[25,0] = 400
$L[26,0] = 900
$L[27,0] = 500
$L[28,0] = 800
$L[29,0] = 600
$L[30,0] = 450
$L[31,0] = 300
$L[32,0] = 800
$L[33,0] = 800
$L[34,0] = 500
$L[35,0] = 350
$L[36,0] = 700
$L[37,0] = 800
$L[38,0] = 700
$L[39,0] = 900
$L[40,0] = 350
$L[41,0] = 600
$L[42,0] = 700
$L[43,0] = 400
$L[44,0] = 500
$L[45,0] = 700
$L[46,0] = 700
$L[47,0] = 450
$L[48,0] = 500
$L[49,0] = 600
$L[50,0] = 500
$L[51,0] = 300
$L[52,0] = 800
$L[53,0] = 800
$L[54,0] = 450
$L[55,0] = 300
$L[56,0] = 700
$L[57,0] = 500
$L[58,0] = 700
$L[59,0] = 700
$L[60,0] = 900
$L[61,0] = 500
$L[62,0] = 350
$L[63,0] = 800
$L[64,0] = 500
$L[65,0] = 350
$L[66,0] = 900
$L[67,0] = 500
$L[68,0] = 350
$L[69,0] = 500
$L[70,0] = 800
$L[71,0] = 700
$L[72,0] = 500
$L[73,0] = 800
$L[74,0] = 1000
$L[75,0] = 450
$L[76,0] = 350
$L[77,0] = 700
$L[78,0] = 400
$L[79,0] = 400
$L[80,0] = 350
$L[81,0] = 800
$L[82,0] = 900
$L[83,0] = 500
$L[84,0] = 900
$L[85,0] = 500
$L[86,0] = 350
$L[87,0] = 500
$L[88,0] = 350
$L[89,0] = 500
$L[90,0] = 400
$L[91,0] = 800
$L[92,0] = 600
$L[93,0] = 800
$L[94,0] = 600
$L[95,0] = 500
$L[96,0] = 500
$L[97,0] = 300
$L[98,0] = 500
$L[99,0] = 800
$L[100,0] = 500
$L[101,0] = 350
$L[102,0] = 500
$L[103,0] = 350
$L[104,0] = 500
$L[105,0] = 400
$L[106,0] = 800
$L[107,0] = 800
$L[108,0] = 600
$L[109,0] = 350
$L[110,0] = 300
$L[111,0] = 500
$L[112,0] = 800
$L[113,0] = 500
$L[114,0] = 500
$L[115,0] = 300
$L[116,0] = 800
$L[117,0] = 350
$L[118,0] = 800
$L[119,0] = 800
$L[120,0] = 600
$L[121,0] = 400
$L[122,0] = 500
$L[123,0] = 500
$L[124,0] = 450
$L[125,0] = 400
$L[126,0] = 400
$L[127,0] = 800
$ws.Range("L88:L215").Value = $L

$M = New-Object 'object[,]' 128,1
$M[0,0] = 600
$M[1,0] = 500
$M[2,0] = 500
$M[3,0] = 800
$M[4,0] = 400
$M[5,0] = 500
$M[6,0] = 400
$M[7,0] = 600
$M[8,0] = 800
$M[9,0] = 500
$M[10,0] = 600
$M[11,0] = 900
$M[12,0] = 800
$M[13,0] = 500
$M[14,0] = 400
$M[15,0] = 500
$M[16,0] = 500
$M[17,0] = 700
$M[18,0] = 450
$M[19,0] = 400
$M[20,0] = 800
$M[21,0] = 400
$M[22,0] = 600
$M[23,0] = 800
$M[24,0] = 400
$M[25,0] = 400
$M[26,0] = 900
$M[27,0] = 500
$M[28,0] = 800
$M[29,0] = 600
$M[30,0] = 450
$M[31,0] = 300
$M[32,0] = 800
$M[33,0] = 800
$M[34,0] = 500
$M[35,0] = 350
$M[36,0] = 700
$M[37,0] = 800
$M[38,0] = 700
$M[39,0] = 900
$M[40,0] = 350
$M[41,0] = 600
$M[42,0] = 700
$M[43,0] = 400
$M[44,0] = 500
$M[45,0] = 700
$M[46,0] = 700
$M[47,0] = 450
$M[48,0] = 500
$M[49,0] = 600
$M[50,0] = 500
$M[51,0] = 300
$M[52,0] = 800
$M[53,0] = 800
$M[54,0] = 450
$M[55,0] = 300
$M[56,0] = 700
$M[57,0] = 500
$M[58,0] = 700
$M[59,0] = 700
$M[60,0] = 900
$M[61,0] = 500
$M[62,0] = 350
$M[63,0] = 800
$M[64,0] = 500
$M[65,0] = 350
$M[66,0] = 900
$M[67,0] = 500
$M[68,0] = 350
$M[69,0] = 500
$M[70,0] = 800
$M[71,0] = 700
$M[72,0] = 500
$M[73,0] = 800
$M[74,0] = 1000
$M[75,0] = 450
$M[76,0] = 350
$M[77,0] = 700
$M[78,0] = 400
$M[79,0] = 400
$M[80,0] = 350
$M[81,0] = 800
$M[82,0] = 900
$M[83,0] = 500
$M[84,0] = 900
$M[85,0] = 500
$M[86,0] = 350
$M[87,0] = 500
$M[88,0] = 350
$M[89,0] = 500
$M[90,0] = 400
$M[91,0] = 800
$M[92,0] = 600
$M[93,0] = 800
$M[94,0] = 600
$M[95,0] = 500
$M[96,0] = 500
$M[97,0] = 300
$M[98,0] = 500
$M[99,0] = 800
$M[100,0] = 500
$M[101,0] = 350
$M[102,0] = 500
$M[103,0] = 350
$M[104,0] = 500
$M[105,0] = 400
$M[106,0] = 800
$M[107,0] = 800
$M[108,0] = 600
$M[109,0] = 350
$M[110,0] = 300
$M[111,0] = 500
$M[112,0] = 800
$M[113,0] = 500
$M[114,0] = 500
$M[115,0] = 300
$M[116,0] = 800
$M[117,0] = 350
$M[118,0] = 800
$M[119,0] = 800
$M[120,0] = 600
$M[121,0] = 400
$M[122,0] = 500
$M[123,0] = 500
$M[124,0] = 450
$M[125,0] = 400
$M[126,0] = 400
$M[127,0] = 800
$ws.Range("M88:M215").Value = $M

$O = New-Object 'object[,]' 128,1
$O[0,0] = 'Región del Maule'
$O[1,0] = 'Región del Maule'
$O[2,0] = 'Región del Maule'
$O[3,0] = 'Región del Maule'
$O[4,0] = 'Provincia de Chacabuco'
$O[5,0] = 'Región del Maule'
$O[6,0] = 'Región del Maule'
$O[7,0] = 'Provincia del Elquí'
$O[8,0] = 'Región del Maule'
$O[9,0] = 'Región del Maule'
$O[10,0] = 'Región del Maule'
$O[11,0] = 'Región del Maule'
$O[12,0] = 'Región del Maule'
$O[13,0] = 'Región del Maule'
$O[14,0] = 'Región del Maule'
$O[15,0] = 'Región del Maule'
$O[16,0] = 'Región del Maule'
$O[17,0] = 'Región del Maule'
$O[18,0] = 'Región del Maule'
$O[19,0] = 'Región del Maule'
$O[20,0] = 'Región del Maule'
$O[21,0] = 'Región del Maule'
$O[22,0] = 'Provincia del Elquí'
$O[23,0] = 'Región del Maule'
$O[24,0] = 'Región del Maule'
$O[25,0] = 'Región del Maule'
$O[26,0] = 'Región del Maule'
$O[27,0] = 'Región del Maule'
$O[28,0] = 'Región del Maule'
$O[29,0] = 'Región del Maule'
$O[30,0] = 'Región del Maule'
$O[31,0] = 'Región del Maule'
$O[32,0] = 'Región del Maule'
$O[33,0] = 'Región del Maule'
$O[34,0] = 'Región del Maule'
$O[35,0] = 'Región del Maule'
$O[36,0] = 'Región del Maule'
$O[37,0] = 'Región del Maule'
$O[38,0] = 'Región del Maule'
$O[39,0] = 'Región del Maule'
$O[40,0] = 'Región del Maule'
$O[41,0] = 'Provincia del Elquí'
$O[42,0] = 'Región del Maule'
$O[43,0] = 'Región del Maule'
$O[44,0] = 'Región del Maule'
$O[45,0] = 'Provincia del Elquí'
$O[46,0] = 'Región del Maule'
$O[47,0] = 'Región del Maule'
$O[48,0] = 'Región del Maule'
$O[49,0] = 'Provincia del Elquí'
$O[50,0] = 'Región del Maule'
$O[51,0] = 'Región del Maule'
$O[52,0] = 'Provincia del Elquí'
$O[53,0] = 'Región del Maule'
$O[54,0] = 'Región del Maule'
$O[55,0] = 'Región del Maule'
$O[56,0] = 'Región del Maule'
$O[57,0] = 'Región del Maule'
$O[58,0] = 'Región del Maule'
$O[59,0] = 'Región del Maule'
$O[60,0] = 'Provincia del Elquí'
$O[61,0] = 'Región del Maule'
$O[62,0] = 'Región del Maule'
$O[63,0] = 'Región del Maule'
$O[64,0] = 'Región del Maule'
$O[65,0] = 'Región del Maule'
$O[66,0] = 'Región del Maule'
$O[67,0] = 'Región del Maule'
$O[68,0] = 'Región del Maule'
$O[69,0] = 'Región del Maule'
$O[70,0] = 'Región del Maule'
$O[71,0] = 'Provincia del Elquí'
$O[72,0] = 'Región del Maule'
$O[73,0] = 'Región del Maule'
$O[74,0] = 'Región del Maule'
$O[75,0] = 'Región del Maule'
$O[76,0] = 'Región del Maule'
$O[77,0] = 'Provincia del Elquí'
$O[78,0] = 'Región del Maule'
$O[79,0] = 'Región del Maule'
$O[80,0] = 'Región del Maule'
$O[81,0] = 'Provincia del Elquí'
$O[82,0] = 'Región del Maule'
$O[83,0] = 'Región del Maule'
$O[84,0] = 'Región del Maule'
$O[85,0] = 'Región del Maule'
$O[86,0] = 'Región del Maule'
$O[87,0] = 'Región del Maule'
$O[88,0] = 'Región del Maule'
$O[89,0] = 'Región del Maule'
$O[90,0] = 'Región del Maule'
$O[91,0] = 'Región del Maule'
$O[92,0] = 'Región del Maule'
$O[93,0] = 'Región del Maule'
$O[94,0] = 'Región del Maule'
$O[95,0] = 'Región del Maule'
$O[96,0] = 'Región del Maule'
$O[97,0] = 'Región del Maule'
$O[98,0] = 'Región del Maule'
$O[99,0] = 'Región del Maule'
$O[100,0] = 'Región del Maule'
$O[101,0] = 'Región del Maule'
$O[102,0] = 'Región del Maule'
$O[103,0] = 'Región del Maule'
$O[104,0] = 'Región del Maule'
$O[105,0] = 'Región del Maule'
$O[106,0] = 'Región del Maule'
$O[107,0] = 'Región del Maule'
$O[108,0] = 'Región del Maule'
$O[109,0] = 'Región del Maule'
$O[110,0] = 'Región del Maule'
$O[111,0] = 'Región del Maule'
$O[112,0] = 'Región del Maule'
$O[113,0] = 'Región del Maule'
$O[114,0] = 'Región del Maule'
$O[115,0] = 'Región del Maule'
$O[116,0] = 'Región del Maule'
$O[117,0] = 'Región del Maule'
$O[118,0] = 'Región del Maule'
$O[119,0] = 'Región del Maule'
$O[120,0] = 'Región del Maule'
$O[121,0] = 'Región del Maule'
$O[122,0] = 'Región del Maule'
$O[123,0] = 'Región del Maule'
$O[124,0] = 'Región del Maule'
$O[125,0] = 'Región del Maule'
$O[126,0] = 'Región del Maule'
$O[127,0] = 'Región del Maule'
$ws.Range("O88:O215").Value = $O

$P = New-Object 'object[,]' 128,1
$P[0,0] = 600
$P[1,0] = 500
$P[2,0] = 500
$P[3,0] = 800
$P[4,0] = 400
$P[5,0] = 500
$P[6,0] = 400
$P[7,0] = 600
$P[8,0] = 800
$P[9,0] = 500
$P[10,0] = 600
$P[11,0] = 900
$P[12,0] = 800
$P[13,0] = 500
$P[14,0] = 400
$P[15,0] = 500
$P[16,0] = 500
$P[17,0] = 700
$P[18,0] = 450
$P[19,0] = 400
$P[20,0] = 800
$P[21,0] = 400
$P[22,0] = 600
$P[23,0] = 800
$P[24,0] = 400
$P[25,0] = 400
$P[26,0] = 900
$P[27,0] = 500
$P[28,0] = 800
$P[29,0] = 600
$P[30,0] = 450
$P[31,0] = 300
$P[32,0] = 800
$P[33,0] = 800
$P[34,0] = 500
$P[35,0] = 350
$P[36,0] = 700
$P[37,0] = 800
$P[38,0] = 700
$P[39,0] = 900
$P[40,0] = 350
$P[41,0] = 600
$P[42,0] = 700
$P[43,0] = 400
$P[44,0] = 500
$P[45,0] = 700
$P[46,0] = 700
$P[47,0] = 450
$P[48,0] = 500
$P[49,0] = 600
$P[50,0] = 500
$P[51,0] = 300
$P[52,0] = 800
$P[53,0] = 800
$P[54,0] = 450
$P[55,0] = 300
$P[56,0] = 700
$P[57,0] = 500
$P[58,0] = 700
$P[59,0] = 700
$P[60,0] = 900
$P[61,0] = 500
$P[62,0] = 350
$P[63,0] = 800
$P[64,0] = 500
$P[65,0] = 350
$P[66,0] = 900
$P[67,0] = 500
$P[68,0] = 350
$P[69,0] = 500
$P[70,0] = 800
$P[71,0] = 700
$P[72,0] = 500
$P[73,0] = 800
$P[74,0] = 1000
$P[75,0] = 450
$P[76,0] = 350
$P[77,0] = 700
$P[78,0] = 400
$P[79,0] = 400
$P[80,0] = 350
$P[81,0] = 800
$P[82,0] = 900
$P[83,0] = 500
$P[84,0] = 900
$P[85,0] = 500
$P[86,0] = 350
$P[87,0] = 500
$P[88,0] = 350
$P[89,0] = 500
$P[90,0] = 400
$P[91,0] = 800
$P[92,0] = 600
$P[93,0] = 800
$P[94,0] = 600
$P[95,0] = 500
$P[96,0] = 500
$P[97,0] = 300
$P[98,0] = 500
$P[99,0] = 800
$P[100,0] = 500
$P[101,0] = 350
$P[102,0] = 500
$P[103,0] = 350
$P[104,0] = 500
$P[105,0] = 400
$P[106,0] = 800
$P[107,0] = 800
$P[108,0] = 600
$P[109,0] = 350
$P[110,0] = 300
$P[111,0] = 500
$P[112,0] = 800
$P[113,0] = 500
$P[114,0] = 500
$P[115,0] = 300
$P[116,0] = 800
$P[117,0] = 350
$P[118,0] = 800
$P[119,0] = 800
$P[120,0] = 600
$P[121,0] = 400
$P[122,0] = 500
$P[123,0] = 500
$P[124,0] = 450
$P[125,0] = 400
$P[126,0] = 400
$P[127,0] = 800
$ws.Range("P88:P215").Value = $P

$ws.Range("D88:D215").NumberFormat = "YYYY-MM-DD HH:MM:SS"